# Map111.xlsx translation patch update:
# - Column B becomes the "translation" column (previously split across B/C/D).
# - For each row, B gets the value that used to live in C (if present) or D
#   (if present); otherwise B just duplicates A's value.
# - Columns C and D are no longer used, so they're cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose B value must come from column C (row -> value).
$fromC = @{
    2 = "Wood"
}

# Rows whose B value must come from column D (row -> value).
$fromD = @{
    10 = "This seems to be the Housemaid's room.`nIt's locked tightly..."
    11 = "The door won't open...`nIt seems to be closed shut by a magical power."
    13 = "Yes"
    14 = "No"
    29 = "This seems to be the Library.`nIt's locked tightly..."
    31 = "This seems to be office of the household owner.`nIt's locked tightly..."
}

for ($r = 1; $r -le 49; $r++) {
    if ($fromC.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $fromC[$r]
    } elseif ($fromD.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $fromD[$r]
    } else {
        $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
    }
}

# Columns C and D are no longer part of the sheet's used range.
$ws.Range("C1:D49").ClearContents()

Write-Output "done"
